$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 30
$ws.Range("F5").Value = 5190
$ws.Range("F6").Value = 5190
$ws.Range("F7").Value = 135
$ws.Range("F8").Value = 134
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = 1170
$ws.Range("F12").Value = 730
$ws.Range("F13").Value = 5105
$ws.Range("F15").Value = 70
$ws.Range("F16").Value = 84
$ws.Range("F17").Value = 246
$ws.Range("F18").Value = 246
$ws.Range("F20").Value = 101
$ws.Range("F22").Value = 3851
$ws.Range("F23").Value = 42
$ws.Range("F24").Value = 3754
$ws.Range("F26").Value = 175
$ws.Range("F28").Value = 226
$ws.Range("F36").Value = 14
$ws.Range("F37").Value = 6673
$ws.Range("F38").Value = 1073
$ws.Range("F43").Value = 1361
$ws.Range("F44").Value = 166
$ws.Range("F45").Value = 673
$ws.Range("F47").Value = 2277
$ws.Range("F49").Value = 91
$ws.Range("F51").Value = 920

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F21").Value = 52
$ws.Range("F24").Value = 810

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 5190
$ws.Range("F8").Value = 5190
$ws.Range("F9").Value = 135
$ws.Range("F10").Value = 134
$ws.Range("F13").Value = 1170
$ws.Range("F14").Value = 730
$ws.Range("F15").Value = 5105
$ws.Range("F17").Value = 70
$ws.Range("F18").Value = 84
$ws.Range("F19").Value = 246
$ws.Range("F20").Value = 247
$ws.Range("F22").Value = 101
$ws.Range("F24").Value = 3852
$ws.Range("F25").Value = 3754
$ws.Range("F27").Value = 175
$ws.Range("F28").Value = 226
$ws.Range("F35").Value = 14
$ws.Range("F37").Value = 6673
$ws.Range("F38").Value = 1073
$ws.Range("F44").Value = 1361
$ws.Range("F45").Value = 166
$ws.Range("F46").Value = 673
$ws.Range("F47").Value = 2277
$ws.Range("F50").Value = 920
